$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("BP4D")
$ws2 = $wb.Worksheets.Item("SEMAINE")

# --- SEMAINE sheet: N1 label text update ---
$ws2.Range("N1").Value2 = "AU28 *might need sep model"

# --- SEMAINE sheet: row 19 values update ---
$ws2.Range("B19").Value2 = 0.16
$ws2.Range("C19").Value2 = 0.25700000000000001
$ws2.Range("D19").Value2 = 0.19719999999999999
$ws2.Range("E19").Value2 = 0.49940000000000001
$ws2.Range("F19").Value2 = 0.43319999999999997
$ws2.Range("G19").Value2 = 0.46379999999999999
$ws2.Range("H19").Value2 = 0.23280000000000001
$ws2.Range("I19").Value2 = 1.8100000000000002E-2
$ws2.Range("J19").Value2 = 3.3399999999999999E-2
$ws2.Range("K19").Value2 = 0.32740000000000002
$ws2.Range("L19").Value2 = 0.24610000000000001
$ws2.Range("M19").Value2 = 0.28089999999999998
$ws2.Range("N19").Value2 = 0
$ws2.Range("O19").Value2 = 0
$ws2.Range("P19").Value2 = 0
$ws2.Range("Q19").Value2 = 0.21940000000000001
$ws2.Range("R19").Value2 = 0.41170000000000001
$ws2.Range("S19").Value2 = 0.2863

# --- BP4D sheet: row 23 values update ---
$ws1.Range("B23").Value2 = 0.38929999999999998
$ws1.Range("C23").Value2 = 0.44629999999999997
$ws1.Range("D23").Value2 = 0.41589999999999999
$ws1.Range("E23").Value2 = 0.3347
$ws1.Range("F23").Value2 = 0.31719999999999998
$ws1.Range("G23").Value2 = 0.32569999999999999
$ws1.Range("H23").Value2 = 0.42399999999999999
$ws1.Range("I23").Value2 = 0.49480000000000002
$ws1.Range("J23").Value2 = 0.45660000000000001
$ws1.Range("K23").Value2 = 0.69499999999999995
$ws1.Range("L23").Value2 = 0.80589999999999995
$ws1.Range("M23").Value2 = 0.74639999999999995
$ws1.Range("N23").Value2 = 0.71740000000000004
$ws1.Range("O23").Value2 = 0.82540000000000002
$ws1.Range("P23").Value2 = 0.76759999999999995
$ws1.Range("Q23").Value2 = 0.8034
$ws1.Range("R23").Value2 = 0.81189999999999996
$ws1.Range("S23").Value2 = 0.80759999999999998
$ws1.Range("T23").Value2 = 0.87109999999999999
$ws1.Range("U23").Value2 = 0.84719999999999995
$ws1.Range("V23").Value2 = 0.85899999999999999
$ws1.Range("W23").Value2 = 0.59589999999999999
$ws1.Range("X23").Value2 = 0.70479999999999998
$ws1.Range("Y23").Value2 = 0.64580000000000004
$ws1.Range("Z23").Value2 = 0.4158
$ws1.Range("AA23").Value2 = 0.5343
$ws1.Range("AB23").Value2 = 0.4677
$ws1.Range("AC23").Value2 = 0.51249999999999996
$ws1.Range("AD23").Value2 = 0.69010000000000005
$ws1.Range("AE23").Value2 = 0.58819999999999995
$ws1.Range("AF23").Value2 = 0.52749999999999997
$ws1.Range("AG23").Value2 = 0.44230000000000003
$ws1.Range("AH23").Value2 = 0.48110000000000003

# --- BP4D sheet: row 27 - rename label, fill values + formulas (was "BP4D MLP geometry", empty) ---
$ws1.Range("A27").Value2 = "BP4D MLP dynamic"

$ws1.Range("B27").Value2 = 0.4219
$ws1.Range("C27").Value2 = 0.44350000000000001
$ws1.Range("D27").Value2 = 0.4325
$ws1.Range("E27").Value2 = 0.38229999999999997
$ws1.Range("F27").Value2 = 0.2429
$ws1.Range("G27").Value2 = 0.29709999999999998
$ws1.Range("H27").Value2 = 0.4874
$ws1.Range("I27").Value2 = 0.48549999999999999
$ws1.Range("J27").Value2 = 0.4864
$ws1.Range("K27").Value2 = 0.76559999999999995
$ws1.Range("L27").Value2 = 0.75790000000000002
$ws1.Range("M27").Value2 = 0.76170000000000004
$ws1.Range("N27").Value2 = 0.72670000000000001
$ws1.Range("O27").Value2 = 0.71279999999999999
$ws1.Range("P27").Value2 = 0.71970000000000001
$ws1.Range("Q27").Value2 = 0.76559999999999995
$ws1.Range("R27").Value2 = 0.85
$ws1.Range("S27").Value2 = 0.80559999999999998
$ws1.Range("T27").Value2 = 0.82969999999999999
$ws1.Range("U27").Value2 = 0.86339999999999995
$ws1.Range("V27").Value2 = 0.84619999999999995
$ws1.Range("W27").Value2 = 0.54079999999999995
$ws1.Range("X27").Value2 = 0.69750000000000001
$ws1.Range("Y27").Value2 = 0.60919999999999996
$ws1.Range("Z27").Value2 = 0.3735
$ws1.Range("AA27").Value2 = 0.34770000000000001
$ws1.Range("AB27").Value2 = 0.36009999999999998
$ws1.Range("AC27").Value2 = 0.59770000000000001
$ws1.Range("AD27").Value2 = 0.47849999999999998
$ws1.Range("AE27").Value2 = 0.53149999999999997
$ws1.Range("AF27").Value2 = 0.4733
$ws1.Range("AG27").Value2 = 0.26669999999999999
$ws1.Range("AH27").Value2 = 0.3412

$ws1.Range("AI27").Formula = "=AVERAGE(B27,E27,H27,K27,N27,Q27,T27,W27,Z27,AC27,AF27)"
$ws1.Range("AJ27").Formula = "=AVERAGE(C27,F27,I27,L27,O27,R27,U27,X27,AA27,AD27,AG27)"
$ws1.Range("AK27").Formula = "=AVERAGE(D27,G27,J27,M27,P27,S27,V27,Y27,AB27,AE27,AH27)"

# --- BP4D sheet: row 28 - clear out label (was "BP4D MLP basic joint scale (working on caelum 104)") ---
$ws1.Range("A28").Value2 = ""
$ws1.Rows.Item(28).RowHeight = 15

# --- Sheet selections / active sheet ---
$ws2.Activate()
$ws2.Range("E30").Select()

$ws1.Activate()
$ws1.Range("AK27").Select()
